# Task: Completed daily operations, 8 hours, 12/01
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new log entry row (row 26) for 12/01/2023
# Copy the date cell's existing number format/style down from the row above
# so the new date cell matches the rest of the column (date-only, no time).
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A26").Value = 45261
$ws.Range("B26").Value = "Internship"
$ws.Range("C26").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

# Update the active selection to reflect the next empty row, like Excel would after data entry
$ws.Range("C27").Select()
